$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '62.869.41'
Set-TextValue 'E2' '  +0.19%  '
Set-TextValue 'D3' '2.463.15'
Set-TextValue 'E3' '  +0.70%  '
Set-TextValue 'E4' '  +0.03%  '
Set-TextValue 'D5' '574.23'
Set-TextValue 'E5' '  -0.32%  '
Set-TextValue 'D6' '146.81'
Set-TextValue 'E6' '  +0.67%  '
Set-TextValue 'E7' '  -0.01%  '
Set-TextValue 'D8' '0.536'
Set-TextValue 'E8' '  -0.57%  '
Set-TextValue 'D9' '2.463.32'
Set-TextValue 'E9' '  +0.77%  '
Set-TextValue 'E10' '  +0.42%  '
Set-TextValue 'D11' '0.163'
Set-TextValue 'E12' '  +0.58%  '
Set-TextValue 'D13' '0.357'
Set-TextValue 'E13' '  +1.11%  '
Set-TextValue 'D14' '29.08'
Set-TextValue 'E14' '  +3.16%  '
Set-TextValue 'E15' '  -0.53%  '
Set-TextValue 'D16' '2.910.58'
Set-TextValue 'E16' '  +0.72%  '
Set-TextValue 'D17' '62.775.61'
Set-TextValue 'E17' '  +0.15%  '
Set-TextValue 'D18' '2.465.71'
Set-TextValue 'E18' '  +0.84%  '
Set-TextValue 'E19' '  +0.34%  '
Set-TextValue 'E20' '  -0.04%  '
Set-TextValue 'D21' '327.10'
Set-TextValue 'E21' '  -0.83%  '
Set-TextValue 'E22' '  +0.01%  '
Set-TextValue 'E23' '  +8.60%  '
Set-TextValue 'D24' '1.00'
Set-TextValue 'E24' '  -0.05%  '
Set-TextValue 'D25' '9.99'
Set-TextValue 'E25' '  +16.89%  '
Set-TextValue 'D26' '65.49'
Set-TextValue 'E26' '  -0.86%  '
Set-TextValue 'D27' '647.80'
Set-TextValue 'E27' '  -0.45%  '
Set-TextValue 'D28' '0.0₃0989'
Set-TextValue 'E28' '  -0.33%  '
Set-TextValue 'E29' '  +0.78%  '
Set-TextValue 'D30' '1.00'
Set-TextValue 'E30' '  -14.95%  '
Set-TextValue 'E31' '  -1.43%  '
Set-TextValue 'E32' '  -2.87%  '
Set-TextValue 'D33' '1.84'
Set-TextValue 'E33' '  -1.26%  '
Set-TextValue 'E34' '  -3.17%  '
Set-TextValue 'D35' '0.998'
Set-TextValue 'E35' '  -0.07%  '
Set-TextValue 'E36' '  +3.00%  '
Set-TextValue 'E37' '  -0.43%  '
Set-TextValue 'E38' '  +4.52%  '
Set-TextValue 'E39' '  -1.49%  '
Set-TextValue 'D40' '5.40'
Set-TextValue 'E40' '  -2.14%  '
Set-TextValue 'D41' '151.43'
Set-TextValue 'E41' '  -1.36%  '
Set-TextValue 'E42' '  -0.34%  '
Set-TextValue 'E43' '  -1.12%  '
Set-TextValue 'E44' '  -48.72%  '
Set-TextValue 'E45' '  -0.01%  '
Set-TextValue 'D46' '152.62'
Set-TextValue 'E46' '  +4.92%  '
Set-TextValue 'D47' '15.26'
Set-TextValue 'E47' '  +2.16%  '
Set-TextValue 'E48' '  -1.63%  '
Set-TextValue 'D49' '20.50'
Set-TextValue 'E49' '  -1.24%  '
Set-TextValue 'E50' '  +0.43%  '
Set-TextValue 'E51' '  -1.34%  '
